$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 76079.75269086983
$ws.Range("D2").Value = 359.5369891933918
$ws.Range("E2").Value = 79427.26180926812
$ws.Range("F2").Value = 375.3566167179011
$ws.Range("I2").Value = 79427.26180926812
$ws.Range("J2").Value = 375.3566167179011
$ws.Range("K2").Value = 19062.54283422435
$ws.Range("L2").Value = 90.08558801229627
$ws.Range("M2").Value = -10931.62729671982
$ws.Range("N2").Value = -51.66058282571956
$ws.Range("Y2").Value = 1536.881987016501
$ws.Range("Z2").Value = 7.26298263090677
$ws.Range("AA2").Value = 196.3929541376
$ws.Range("AB2").Value = 0.928112
$ws.Range("AE2").Value = 79427.26180926812
$ws.Range("AF2").Value = 375.3566167179011
$ws.Range("AM2").Value = 23344.91068630396
$ws.Range("AN2").Value = 110.3231622642963
$ws.Range("C3").Value = 228669.6487221458
$ws.Range("D3").Value = 331.2785379538058
$ws.Range("E3").Value = 229126.9880195901
$ws.Range("F3").Value = 331.9410950297134
$ws.Range("I3").Value = 229126.9880195901
$ws.Range("J3").Value = 331.9410950297134
$ws.Range("K3").Value = 41242.85784352622
$ws.Range("L3").Value = 59.74939710534842
$ws.Range("M3").Value = -19068.40789000158
$ws.Range("N3").Value = -27.62480426329866
$ws.Range("Y3").Value = 11957.81677241365
$ws.Range("Z3").Value = 17.32354109791855
$ws.Range("AA3").Value = 1465.38644026243
$ws.Range("AB3").Value = 2.122936210294076
$ws.Range("AE3").Value = 229126.9880195901
$ws.Range("AF3").Value = 331.9410950297134
$ws.Range("AM3").Value = 52356.78303622502
$ws.Range("AN3").Value = 75.85037469174843
$ws.Range("C4").Value = 54035.40589132848
$ws.Range("D4").Value = 591.1122792174041
$ws.Range("E4").Value = 55548.39725628567
$ws.Range("F4").Value = 607.6634230354913
$ws.Range("I4").Value = 55548.39725628567
$ws.Range("J4").Value = 607.6634230354913
$ws.Range("K4").Value = 17220.00314944856
$ws.Range("L4").Value = 188.3756611410023
$ws.Range("M4").Value = -16731.51865307658
$ws.Range("N4").Value = -183.0319577071184
$ws.Range("Y4").Value = 920.3545305074829
$ws.Range("Z4").Value = 10.06808138557256
$ws.Range("AA4").Value = 70.8326679330217
$ws.Range("AB4").Value = 0.774863426937952
$ws.Range("AE4").Value = 55548.39725628567
$ws.Range("AF4").Value = 607.6634230354913
$ws.Range("AM4").Value = 48865.06795812729
$ws.Range("AN4").Value = 534.5521370364563
$ws.Range("C5").Value = 246966.0511828376
$ws.Range("D5").Value = 371.5758247263654
$ws.Range("E5").Value = 253881.100615957
$ws.Range("F5").Value = 381.9799478187036
$ws.Range("I5").Value = 253881.100615957
$ws.Range("J5").Value = 381.9799478187036
$ws.Range("K5").Value = 78703.14119094668
$ws.Range("L5").Value = 118.4137838237981
$ws.Range("M5").Value = -19849.72800677411
$ws.Range("N5").Value = -29.86515360870653
$ws.Range("Y5").Value = 53488.97833993503
$ws.Range("Z5").Value = 80.47750346754235
$ws.Range("AA5").Value = 11176.73858566045
$ws.Range("AB5").Value = 16.8161001798711
$ws.Range("AE5").Value = 253881.100615957
$ws.Range("AF5").Value = 381.9799478187036
$ws.Range("AM5").Value = 106061.3949179582
$ws.Range("AN5").Value = 159.5759826077981
$ws.Range("AP5").Value = 51157.71129149164
$ws.Range("AQ5").Value = 24093.34654821123
$ws.Range("AR5").Value = 4911.914497344042
$ws.Range("AS5").Value = 3423.393550588788
$ws.Range("AT5").Value = 18729.05669534759
$ws.Range("C6").Value = 97438.4374241778
$ws.Range("D6").Value = 167.6485968223643
$ws.Range("E6").Value = 102700.1130450834
$ws.Range("F6").Value = 176.701621050772
$ws.Range("I6").Value = 102700.1130450834
$ws.Range("J6").Value = 176.701621050772
$ws.Range("K6").Value = 24648.02713082001
$ws.Range("L6").Value = 42.40838905218527
$ws.Range("M6").Value = -32124.4185757937
$ws.Range("N6").Value = -55.27196289613071
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 0
$ws.Range("AA6").Value = 13904.67036138248
$ws.Range("AB6").Value = 23.92380806780943
$ws.Range("AE6").Value = 102700.1130450834
$ws.Range("AF6").Value = 176.701621050772
$ws.Range("AM6").Value = 42983.77467323736
$ws.Range("AN6").Value = 73.95612758778526
$ws.Range("C7").Value = 85827.38362134235
$ws.Range("D7").Value = 119.249583692393
$ws.Range("E7").Value = 88230.55036273993
$ws.Range("F7").Value = 122.58857203578
$ws.Range("I7").Value = 88230.55036273993
$ws.Range("J7").Value = 122.58857203578
$ws.Range("K7").Value = 27351.47061244938
$ws.Range("L7").Value = 38.00245733109182
$ws.Range("M7").Value = -234405.860343976
$ws.Range("N7").Value = -325.6862796191011
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 0
$ws.Range("Y7").Value = 0
$ws.Range("Z7").Value = 0
$ws.Range("AA7").Value = 17960.1992167857
$ws.Range("AB7").Value = 24.95411358551024
$ws.Range("AE7").Value = 88230.55036273993
$ws.Range("AF7").Value = 122.58857203578
$ws.Range("AM7").Value = 190584.0058627335
$ws.Range("AN7").Value = 264.7996757984374
$ws.Range("C8").Value = 75389.14597457208
$ws.Range("D8").Value = 111.0695995048766
$ws.Range("E8").Value = 78706.26839745326
$ws.Range("F8").Value = 115.9566618830912
$ws.Range("I8").Value = 78706.26839745326
$ws.Range("J8").Value = 115.9566618830912
$ws.Range("K8").Value = 18889.50441538878
$ws.Range("L8").Value = 27.82959885194189
$ws.Range("M8").Value = -68510.2887758827
$ws.Range("N8").Value = -100.9350913574124
$ws.Range("Y8").Value = 16641.31830619365
$ws.Range("Z8").Value = 24.51738291511521
$ws.Range("AA8").Value = 1226.75045097911
$ws.Range("AB8").Value = 1.807351436619718
$ws.Range("AE8").Value = 78706.26839745326
$ws.Range("AF8").Value = 115.9566618830912
$ws.Range("AM8").Value = 60102.78762797349
$ws.Range("AN8").Value = 88.54845700490189
$ws.Range("C9").Value = 52221.25614183004
$ws.Range("D9").Value = 207.9235031194944
$ws.Range("E9").Value = 52325.6986541137
$ws.Range("F9").Value = 208.3393501257333
$ws.Range("I9").Value = 52325.6986541137
$ws.Range("J9").Value = 208.3393501257333
$ws.Range("K9").Value = 9418.625757740467
$ws.Range("L9").Value = 37.50108302263201
$ws.Range("M9").Value = -7377.168810520629
$ws.Range("N9").Value = -29.37284346476406
$ws.Range("Y9").Value = 1784.70949962666
$ws.Range("Z9").Value = 7.10597711792252
$ws.Range("AA9").Value = 135.3772805322241
$ws.Range("AB9").Value = 0.5390164942528733
$ws.Range("AE9").Value = 52325.6986541137
$ws.Range("AF9").Value = 208.3393501257333
$ws.Range("AM9").Value = 14274.31596622362
$ws.Range("AN9").Value = 56.83443868663201
$ws.Range("AP9").Value = 20009.02046190888
$ws.Range("AQ9").Value = 8211.543904664404
$ws.Range("AR9").Value = 3324.000478051049
$ws.Range("AS9").Value = 2253.505102732599
$ws.Range("AT9").Value = 6219.970976460829
$ws.Range("C10").Value = 124413.9281094219
$ws.Range("D10").Value = 192.7382500294682
$ws.Range("E10").Value = 129888.1409462365
$ws.Range("F10").Value = 201.2187330307648
$ws.Range("I10").Value = 129888.1409462365
$ws.Range("J10").Value = 201.2187330307648
$ws.Range("K10").Value = 31173.15382709676
$ws.Range("L10").Value = 48.29249592738356
$ws.Range("M10").Value = -28428.19120649435
$ws.Range("N10").Value = -44.04008383871528
$ws.Range("Y10").Value = 5301.670831492414
$ws.Range("Z10").Value = 8.213186206896552
$ws.Range("AA10").Value = 416.0577927167999
$ws.Range("AB10").Value = 0.6445439999999998
$ws.Range("AE10").Value = 129888.1409462365
$ws.Range("AF10").Value = 201.2187330307648
$ws.Range("AM10").Value = 57395.70934482056
$ws.Range("AN10").Value = 88.91567645538355
$ws.Range("AP10").Value = 56884.77943034428
$ws.Range("AQ10").Value = 26722.02729027598
$ws.Range("AR10").Value = 5820.37064949248
$ws.Range("AS10").Value = 3639.313128617208
$ws.Range("AT10").Value = 20703.06836195861
